$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '64.103.30'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  +2.14%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.101.39'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  +1.64%  '
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '545.39'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +0.03%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '141.25'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +5.21%  '
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.098.30'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +1.75%  '
$ws.Range("E9").Value = '  +2.38%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.63'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +3.90%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.157'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +1.72%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.460'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +0.46%  '
$ws.Range("E13").Value = '  +6.40%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '35.05'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +1.01%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.612.60'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +1.85%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '64.233.85'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +2.30%  '
$ws.Range("E17").Value = '  +1.88%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.107.30'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +1.84%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.70'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +1.28%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '484.62'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +0.63%  '
$ws.Range("E21").Value = '  +1.28%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.703'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +0.95%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.16'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +1.80%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '79.49'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +3.03%  '
$ws.Range("E25").Value = '  +2.69%  '
$ws.Range("E26").Value = '  +0.14%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.74'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +1.65%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.14'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -1.67%  '
$ws.Range("E29").Value = '  -0.05%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '26.44'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +1.31%  '
$ws.Range("E31").Value = '  -0.77%  '
$ws.Range("E32").Value = '  +3.55%  '
$ws.Range("B33").Value = 'Stacks'
$ws.Range("C33").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.38'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -3.84%  '
$ws.Range("B34").Value = 'OKB'
$ws.Range("C34").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '57.71'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -0.99%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.43'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +6.96%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '499.08'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -2.38%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.04'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +1.70%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.293.89'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +7.01%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0409'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +2.93%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0804'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +2.35%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.121'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +2.63%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.76'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +6.12%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.15'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +1.75%  '
$ws.Range("E44").Value = '  +1.42%  '
$ws.Range("B46").Value = 'InjectiveProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '25.47'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +4.62%  '
$ws.Range("B47").Value = 'Monero'
$ws.Range("C47").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '123.99'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +3.40%  '
$ws.Range("B48").Value = 'Fetch.AI'
$ws.Range("C48").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.06'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +1.80%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0₃0537'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +8.97%  '
$ws.Range("E50").Value = '  +2.82%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.40'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +0.19%  '
